# CA Workload Automation v25.pptx - fix casing of the "summary" bullet on
# the "technical agenda" slide (slide 2, body placeholder "Rectangle 3").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$hit = $textRange.Find("summary", 0)
if ($hit -ne $null) {
    $hit.Text = "Summary"
}
